$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# For column D, temporarily force Text format so values such as "1.003" or
# "0.06286" are stored as text (matching the original inlineStr cells) rather
# than being auto-converted to numbers by Excel, then restore the original style.

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.209.06'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  -0.76%  '
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.671.93'
$ws.Range('D3').Style = $origStyle
$ws.Range('E3').Value = '  -1.43%  '
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  -0.69%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.05'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -2.91%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5283'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -3.60%  '
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('E8').Value = '  -3.32%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06286'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  -2.45%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.30'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -3.04%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07555'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -1.61%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.656.61'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  -2.61%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.458'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  -2.17%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5598'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -4.22%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.14'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +2.30%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008031'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  -4.48%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.263.10'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('E18').Value = '  -0.71%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.795'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -3.00%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '187.42'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -2.08%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.41'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -5.26%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.209'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  -0.89%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -0.65%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.96'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +0.75%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1260'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -3.95%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.581'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -4.27%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.97'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +1.12%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06203'
$ws.Range('D28').Style = $origStyle
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.366'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -1.60%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.285'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -3.36%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.507'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  -2.83%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.431'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -4.54%  '
$ws.Range('E33').Value = '  -3.09%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.001'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  -3.60%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6074'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range('E36').Value = '  -0.03%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.745'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  -0.58%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.123'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +0.34%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01624'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -1.87%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.101.20'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -1.54%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8776'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('E42').Value = '  -0.98%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.99'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('E44').Value = '  -1.29%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000107'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -1.87%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.04'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  -2.74%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.007'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +0.03%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.057'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  -1.69%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05224'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('E50').Value = '  -1.13%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.993'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -2.14%  '
